# Add a new "isNewPlayer" column (D) to the PlayerInfo sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rows (rows 1-3 mirror the structure of columns A-C)
$ws.Range("D1").Value = "isNewPlayer"
$ws.Range("D2").Value = "bool"
$ws.Range("D3").Value = "is new player"

# Data rows: first two entered as boolean TRUE/FALSE, rest as plain 0/1 numbers
$ws.Range("D4").Value = $true
$ws.Range("D5").Value = $false
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1

$ws.Range("E11").Select()
